$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.227462887763977
$ws.Range("B1").Value = 2.795978546142578
$ws.Range("C1").Value = 3.417288541793823
$ws.Range("D1").Value = 3.80321216583252
$ws.Range("E1").Value = 0.8894218802452087
